# Auto-generated edit script for game_06 skill tree redesign
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A, B, C, F contain zero-padded / small-integer strings (e.g. "01", "0001")
# that Excel would otherwise auto-convert to numbers, stripping leading zeros.
# Force those columns to Text format for the full data range before writing values.
$ws.Range("A6:C17").NumberFormat = "@"
$ws.Range("F6:F17").NumberFormat = "@"

# Row 6
$ws.Range('A6').Value = '70'
$ws.Range('B6').Value = '01'
$ws.Range('C6').Value = '0001'
$ws.Range('D6').Value = '裂变弹匣'
$ws.Range('E6').Value = '弹道'
$ws.Range('F6').Value = '1'
$ws.Range('G6').Value = ''
$ws.Range('H6').Value = 'split:+1|splitAngle:+12|damage:+6'
$ws.Range('I6').Value = 'level:3'
$ws.Range('J6').Value = '将弹道裂变为额外子弹，基础火力提升。'
$ws.Range('K6').Value = 'icons/skill/focal_anchor.png'

# Row 7
$ws.Range('A7').Value = '70'
$ws.Range('B7').Value = '01'
$ws.Range('C7').Value = '0002'
$ws.Range('D7').Value = '轨迹稳流'
$ws.Range('E7').Value = '弹道'
$ws.Range('F7').Value = '2'
$ws.Range('G7').Value = 'skill:70010001'
$ws.Range('H7').Value = 'stability:+14|projectileSize:+18|projectileSpeed:+8'
$ws.Range('I7').Value = 'level:6'
$ws.Range('J7').Value = '导流装置压制散布，并扩大弹道厚度。'
$ws.Range('K7').Value = 'icons/skill/rapid_siphon.png'

# Row 8
$ws.Range('A8').Value = '70'
$ws.Range('B8').Value = '01'
$ws.Range('C8').Value = '0003'
$ws.Range('D8').Value = '深域贯穿'
$ws.Range('E8').Value = '弹道'
$ws.Range('F8').Value = '3'
$ws.Range('G8').Value = 'skill:70010002'
$ws.Range('H8').Value = 'pierce:+2|crit:+5|damage:+10'
$ws.Range('I8').Value = 'level:9'
$ws.Range('J8').Value = '强化穿甲结构并提升暴击输出。'
$ws.Range('K8').Value = 'icons/skill/zero_point.png'

# Row 9
$ws.Range('A9').Value = '70'
$ws.Range('B9').Value = '02'
$ws.Range('C9').Value = '0001'
$ws.Range('D9').Value = '棱镜导光'
$ws.Range('E9').Value = '能量'
$ws.Range('F9').Value = '1'
$ws.Range('G9').Value = ''
$ws.Range('H9').Value = 'damage:+6|sanityDrain:-6|projectileSize:+18'
$ws.Range('I9').Value = 'level:4|weaponAttack:BEAM'
$ws.Range('J9').Value = '棱镜束缚能量消耗，同时扩大光束宽度。'
$ws.Range('K9').Value = 'icons/skill/aether_resonance.png'

# Row 10
$ws.Range('A10').Value = '70'
$ws.Range('B10').Value = '02'
$ws.Range('C10').Value = '0002'
$ws.Range('D10').Value = '谐振折叠'
$ws.Range('E10').Value = '能量'
$ws.Range('F10').Value = '2'
$ws.Range('G10').Value = 'skill:70020001'
$ws.Range('H10').Value = 'damageMultiplier:+8|projectileSpeed:+14|stability:+8'
$ws.Range('I10').Value = 'level:7|weaponAttack:BEAM'
$ws.Range('J10').Value = '折叠振镜提高能量聚焦与射速。'
$ws.Range('K10').Value = 'icons/skill/choir_surge.png'

# Row 11
$ws.Range('A11').Value = '70'
$ws.Range('B11').Value = '02'
$ws.Range('C11').Value = '0003'
$ws.Range('D11').Value = '相干放射'
$ws.Range('E11').Value = '能量'
$ws.Range('F11').Value = '3'
$ws.Range('G11').Value = 'skill:70020002'
$ws.Range('H11').Value = 'ricochet:+1|crit:+6|damage:+12'
$ws.Range('I11').Value = 'level:10|weaponAttack:BEAM'
$ws.Range('J11').Value = '相干腔反复震荡，使光束可在敌间折射。'
$ws.Range('K11').Value = 'icons/skill/zero_point.png'

# Row 12
$ws.Range('A12').Value = '70'
$ws.Range('B12').Value = '03'
$ws.Range('C12').Value = '0001'
$ws.Range('D12').Value = '相位壁垒'
$ws.Range('E12').Value = '护卫'
$ws.Range('F12').Value = '1'
$ws.Range('G12').Value = ''
$ws.Range('H12').Value = 'shield:+60|contactResist:+25|sanityRegen:+3'
$ws.Range('I12').Value = 'level:4'
$ws.Range('J12').Value = '展开相位护壁，降低接触伤害并补充理智。'
$ws.Range('K12').Value = 'icons/skill/ward_bastion.png'

# Row 13
$ws.Range('A13').Value = '70'
$ws.Range('B13').Value = '03'
$ws.Range('C13').Value = '0002'
$ws.Range('D13').Value = '护盾崩击'
$ws.Range('E13').Value = '护卫'
$ws.Range('F13').Value = '2'
$ws.Range('G13').Value = 'skill:70030001'
$ws.Range('H13').Value = 'meleeDamage:+70|meleeRadius:+20|meleeInterval:-0.5'
$ws.Range('I13').Value = 'level:7'
$ws.Range('J13').Value = '护盾冲击形成短距爆发，持续清理近身威胁。'
$ws.Range('K13').Value = 'icons/skill/seraphic_shell.png'

# Row 14
$ws.Range('A14').Value = '70'
$ws.Range('B14').Value = '03'
$ws.Range('C14').Value = '0003'
$ws.Range('D14').Value = '寂光回响'
$ws.Range('E14').Value = '护卫'
$ws.Range('F14').Value = '3'
$ws.Range('G14').Value = 'skill:70030002'
$ws.Range('H14').Value = 'beamReflect:20%|shieldRegen:+16|invulnTime:+0.4'
$ws.Range('I14').Value = 'level:10'
$ws.Range('J14').Value = '护盾折射寂光，可短暂反弹能量。'
$ws.Range('K14').Value = 'icons/skill/ward_bastion.png'

# Row 15
$ws.Range('A15').Value = '70'
$ws.Range('B15').Value = '04'
$ws.Range('C15').Value = '0001'
$ws.Range('D15').Value = '术式镀层'
$ws.Range('E15').Value = '工坊'
$ws.Range('F15').Value = '1'
$ws.Range('G15').Value = ''
$ws.Range('H15').Value = 'projectileSize:+24|elementSlow:+18|elementSlowDuration:+1.4'
$ws.Range('I15').Value = 'level:5'
$ws.Range('J15').Value = '在弹体上刻蚀术式，对命中目标施加霜蚀减速。'
$ws.Range('K15').Value = 'icons/skill/undertow.png'

# Row 16
$ws.Range('A16').Value = '70'
$ws.Range('B16').Value = '04'
$ws.Range('C16').Value = '0002'
$ws.Range('D16').Value = '弹道精铸'
$ws.Range('E16').Value = '工坊'
$ws.Range('F16').Value = '2'
$ws.Range('G16').Value = 'skill:70040001'
$ws.Range('H16').Value = 'split:+1|splitAngle:+4|pierce:+1'
$ws.Range('I16').Value = 'level:8'
$ws.Range('J16').Value = '精铸枪管令术弹再次分裂并保持贯穿。'
$ws.Range('K16').Value = 'icons/skill/riptide_collapse.png'

# Row 17
$ws.Range('A17').Value = '70'
$ws.Range('B17').Value = '04'
$ws.Range('C17').Value = '0003'
$ws.Range('D17').Value = '秘火迸流'
$ws.Range('E17').Value = '工坊'
$ws.Range('F17').Value = '3'
$ws.Range('G17').Value = 'skill:70040002'
$ws.Range('H17').Value = 'damageMultiplier:+12|luckBonus:+12|projectileSpeed:+16'
$ws.Range('I17').Value = 'level:11'
$ws.Range('J17').Value = '秘火符文强化弹速与掉落运势。'
$ws.Range('K17').Value = 'icons/skill/choir_surge.png'

# Update the ignoredErrors / dimension sqref implicitly handled by Excel on save.